$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (Receptor avg/total expression and derived specificities recalculated with new TPM)
$ws.Range("M2").Value = 0.06166766666666667
$ws.Range("N2").Value = 0.185003
$ws.Range("O2").Value = 0.3189772891852935
$ws.Range("P2").Value = 0.3189772891852935
$ws.Range("Q2").Value = 0.0001012171968888889
$ws.Range("R2").Value = 0.0009109547720000001
$ws.Range("S2").Value = 0.3189772891852935
$ws.Range("T2").Value = 0.3189772891852935

# Row 3 updates (derived specificities recalculated with new TPM)
$ws.Range("O3").Value = 0.4045463009579509
$ws.Range("P3").Value = 0.4045463009579509
$ws.Range("S3").Value = 0.4045463009579509
$ws.Range("T3").Value = 0.4045463009579509

# Row 4 updates (Receptor-expressing cells, detection rate, avg/total expression and derived specificities)
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.05345100000000001
$ws.Range("N4").Value = 0.160353
$ws.Range("O4").Value = 0.2764764098567557
$ws.Range("P4").Value = 0.2764764098567557
$ws.Range("Q4").Value = 0.00008773090800000002
$ws.Range("R4").Value = 0.0007895781720000002
$ws.Range("S4").Value = 0.2764764098567557
$ws.Range("T4").Value = 0.2764764098567557
